$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.187.74'
$ws.Range('E2').Value = '  -4.22%  '

# Row 3
$ws.Range('D3').Value = '3.537.74'
$ws.Range('E3').Value = '  -5.04%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.32'
$ws.Range('E5').Value = '  -7.35%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.47'
$ws.Range('E6').Value = '  -2.92%  '

# Row 7
$ws.Range('D7').Value = '3.533.41'
$ws.Range('E7').Value = '  -5.00%  '

# Row 8
$ws.Range('E8').Value = '  -5.11%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.01%  '

# Row 10
$ws.Range('E10').Value = '  -9.21%  '

# Row 11
$ws.Range('E11').Value = '  -10.84%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.88'
$ws.Range('E12').Value = '  -10.75%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  -12.89%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.68'
$ws.Range('E14').Value = '  -8.28%  '

# Row 15
$ws.Range('D15').Value = '4.097.40'
$ws.Range('E15').Value = '  -5.36%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.523.99'
$ws.Range('E16').Value = '  -5.65%  '

# Row 17
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.125'
$ws.Range('E17').Value = '  -1.55%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '66.222.08'
$ws.Range('E18').Value = '  -4.05%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.08'
$ws.Range('E19').Value = '  -7.60%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.99'
$ws.Range('E20').Value = '  -7.86%  '

# Row 21
$ws.Range('E21').Value = '  -9.37%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '385.72'
$ws.Range('E22').Value = '  -6.69%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.17'
$ws.Range('E23').Value = '  -9.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.81'
$ws.Range('E24').Value = '  -6.19%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.98'
$ws.Range('E25').Value = '  -4.18%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.90'
$ws.Range('E26').Value = '  -6.63%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.29'
$ws.Range('E27').Value = '  -5.58%  '

# Row 28
$ws.Range('E28').Value = '  -0.16%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.52'
$ws.Range('E29').Value = '  -8.02%  '

# Row 30
$ws.Range('E30').Value = '  -10.66%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.46'
$ws.Range('E31').Value = '  -5.35%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.75'
$ws.Range('E32').Value = '  -6.86%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.13'
$ws.Range('E33').Value = '  -5.28%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '616.73'
$ws.Range('E34').Value = '  -5.10%  '

# Row 35
$ws.Range('E35').Value = '  -8.20%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '62.78'
$ws.Range('E36').Value = '  -6.34%  '

# Row 37
$ws.Range('E37').Value = '  -10.61%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.400'
$ws.Range('E38').Value = '  -4.12%  '

# Row 39
$ws.Range('E39').Value = '  +0.17%  '

# Row 40
$ws.Range('D40').Value = '0.0₃0728'
$ws.Range('E40').Value = '  -13.13%  '

# Row 41
$ws.Range('E41').Value = '  -6.74%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.16%  '

# Row 43
$ws.Range('D43').Value = '3.068.47'
$ws.Range('E43').Value = '  +5.46%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.89'
$ws.Range('E44').Value = '  -6.39%  '

# Row 45
$ws.Range('E45').Value = '  -4.65%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0403'
$ws.Range('E46').Value = '  -10.58%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.16'
$ws.Range('E47').Value = '  +2.10%  '

# Row 48
$ws.Range('E48').Value = '  -8.42%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '137.23'
$ws.Range('E49').Value = '  -4.40%  '

# Row 50
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.35'
$ws.Range('E50').Value = '  -9.58%  '

# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.72'
$ws.Range('E51').Value = '  -2.56%  '
